# ImportPeople_template.xlsx edit:
# - store IdentityNumber (col C) and ManagerId (col K) as text instead of numbers
# - add two new data rows (6 and 7) so several contracts can be grouped
#   under one account (identity id / manager id)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header style tweaks: C1 and K1 move to the new "text" style (numFmtId 49)
# ---------------------------------------------------------------------------
$ws.Range("C1").NumberFormat = "@"
$ws.Range("K1").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Existing IdentityNumber / ManagerId values become text instead of numbers.
# Order matters here: it controls the order new shared strings are created in.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "AOUDI1"
$ws.Range("A7").Value = "AOUDI1"
$ws.Range("B6").Value = "JIHENE1"
$ws.Range("B7").Value = "JIHENE1"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "290109933804360"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "290109933804360"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "168039933804012"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "169129934203490"

$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "290109933804360"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "168039933804012"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1680399338040120"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1680399338040120"

$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "290109933804360"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "169129934203490"

# ---------------------------------------------------------------------------
# Row 6 - new account/contract (copy of row 4's data, grouped under AOUDI1)
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = 1303026

$ws.Range("E6").Value = "othermail@mail.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:othermail@mail.com") | Out-Null
$ws.Range("E6").Style = "Normal"

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").HorizontalAlignment = -4152
$ws.Range("F6").Value = "+79213456791"

$ws.Range("G6").NumberFormat = "dd/mm/yy;@"
$ws.Range("G6").Value = 36537

$ws.Range("H6").Value = "RU"
$ws.Range("I6").Value = "SK"
$ws.Range("J6").Value = "Female"

$ws.Range("L6").Value = 1

$ws.Range("M6").VerticalAlignment = -4160
$ws.Range("M6").Value = "ACCOUNTING JUNIOR MANAGER"

$ws.Range("N6").Value = 1

$ws.Range("O6").NumberFormat = "m/d/yy"
$ws.Range("O6").Value = 42736
$ws.Range("P6").NumberFormat = "m/d/yy"

$ws.Range("Q6").Value = "Active"

# ---------------------------------------------------------------------------
# Row 7 - second contract for the same account (AOUDI1 / JIHENE1)
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = 1303026

$ws.Range("E7").Value = "othermail@mail.com"
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:othermail@mail.com") | Out-Null
$ws.Range("E7").Style = "Normal"

$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("F7").Value = "+79213456791"

$ws.Range("G7").NumberFormat = "dd/mm/yy;@"
$ws.Range("G7").Value = 36537

$ws.Range("H7").Value = "RU"
$ws.Range("I7").Value = "SK"
$ws.Range("J7").Value = "Female"

$ws.Range("L7").Value = 2

$ws.Range("M7").VerticalAlignment = -4160
$ws.Range("M7").Value = "ACCOUNTING JUNIOR MANAGER"

$ws.Range("N7").Value = 1

$ws.Range("O7").NumberFormat = "m/d/yy"
$ws.Range("O7").Value = 42736
$ws.Range("P7").NumberFormat = "m/d/yy"

$ws.Range("Q7").Value = "Active"

# ---------------------------------------------------------------------------
# Selection moved by the author while reviewing the new rows
# ---------------------------------------------------------------------------
$ws.Range("I3").Select() | Out-Null

# Remove the auto-generated "Hyperlink" cell style left behind by
# Hyperlinks.Add - the source file never used Excel's hyperlink theming.
$wb.Styles.Item("Hyperlink").Delete()

Write-Host "ImportPeople_template.xlsx updated"
